# Generate Report for Handback
# Updates the localization-status workbook: marks the de-de/zh-cn handback
# rows as complete (target + handback file names, handback datetime) and
# updates the Overview "Status" column text, widening columns so the new,
# longer text/hyperlinks fit.

$wb = $excel.ActiveWorkbook

$mdUrl04 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8c2fa64e19c3ecf18db8cb48e9582b4782dcb279/e2e/04f690e6-e310-4e91-83c7-d9485eadbb02.md"
$mdUrl6b = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8c2fa64e19c3ecf18db8cb48e9582b4782dcb279/e2e/6bcce489-181b-4ee7-b582-ce837ab9b595.md"

function Set-HyperlinkCell($ws, $cellRef, $url, $text) {
    $ws.Hyperlinks.Add($ws.Range($cellRef), $url, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $text)
    $rng = $ws.Range($cellRef)
    $rng.Style = "HyperLink"
    $rng.Font.Underline = 2
    # COLORREF (BGR) encoding of RGB(100,149,237) = #6495ED -- matches the
    # workbook's existing custom "HyperLink" cell style font color.
    $rng.Font.Color = 15570276
}

# ---------------------------------------------------------------------
# Overview sheet: handback status text for both locale rows
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# Widen the status columns (E, F) so the longer text is visible.
$wsOverview.Range("E1").ColumnWidth = 29.9777047293527 - (5/6)
$wsOverview.Range("F1").ColumnWidth = 29.9777047293527 - (5/6)

# ---------------------------------------------------------------------
# zh-cn sheet: populate Latest Target File / Latest Handback File /
# Latest Handback DateTime for both rows
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

Set-HyperlinkCell $wsZhCn "I2" $mdUrl04 "04f690e6-e310-4e91-83c7-d9485eadbb02.md"
$wsZhCn.Range("J2").Value = "04f690e6-e310-4e91-83c7-d9485eadbb02.f6624f78fd80e0b49a692d9b95c4c576099b97ba.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-09-05 15:00:16"

Set-HyperlinkCell $wsZhCn "I3" $mdUrl6b "6bcce489-181b-4ee7-b582-ce837ab9b595.md"
$wsZhCn.Range("J3").Value = "6bcce489-181b-4ee7-b582-ce837ab9b595.1a2690f9f86538d76be41880767fe3aefcb9144b.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-09-05 15:00:16"

$wsZhCn.Range("C1").ColumnWidth = 29.9777047293527 - (5/6)
$wsZhCn.Range("I1").ColumnWidth = 40 - (5/6)
$wsZhCn.Range("J1").ColumnWidth = 40 - (5/6)

# ---------------------------------------------------------------------
# de-de sheet: populate Latest Target File / Latest Handback File /
# Latest Handback DateTime for both rows
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

Set-HyperlinkCell $wsDeDe "I2" $mdUrl04 "04f690e6-e310-4e91-83c7-d9485eadbb02.md"
$wsDeDe.Range("J2").Value = "04f690e6-e310-4e91-83c7-d9485eadbb02.f6624f78fd80e0b49a692d9b95c4c576099b97ba.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-09-05 15:00:44"

Set-HyperlinkCell $wsDeDe "I3" $mdUrl6b "6bcce489-181b-4ee7-b582-ce837ab9b595.md"
$wsDeDe.Range("J3").Value = "6bcce489-181b-4ee7-b582-ce837ab9b595.1a2690f9f86538d76be41880767fe3aefcb9144b.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-09-05 15:00:44"

$wsDeDe.Range("C1").ColumnWidth = 29.9777047293527 - (5/6)
$wsDeDe.Range("I1").ColumnWidth = 40 - (5/6)
$wsDeDe.Range("J1").ColumnWidth = 40 - (5/6)

Write-Host "Edit complete"
